$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: merge the "To-do" + bookmark + ":  " runs in paragraph 1 into a
# single run "To-do:  ", dropping the _GoBack bookmark from that location.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$firstPara = $d.Paragraphs.Item(1)
$firstRange = $firstPara.Range
$firstRange.MoveEnd(1, -1) | Out-Null       # exclude the paragraph mark
$firstRange.Delete()
$firstRange2 = $d.Paragraphs.Item(1).Range
$firstRange2.MoveEnd(1, -1) | Out-Null
$firstRange2.InsertAfter("To-do:  ")

# ---------------------------------------------------------------------------
# Change 2: after the "delete characters" bullet item, append two new bullet
# paragraphs (same ListParagraph / numId 1 style) and re-anchor the _GoBack
# bookmark at the end of the last of these new paragraphs.
# ---------------------------------------------------------------------------
$anchorText = "As a user I want to be able to delete characters that I don" + [char]0x2019 + "t want anymore "

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $paraText = $d.Paragraphs.Item($i).Range.Text
    if ($paraText.TrimEnd([char]13) -eq $anchorText) {
        $targetPara = $d.Paragraphs.Item($i)
        break
    }
}

$newPara1 = $targetPara.Range.InsertParagraphAfter()
$para1Index = $targetPara.Index + 1
$d.Paragraphs.Item($para1Index).Range.InsertBefore("As a user I want to be able to search through a list of spells to select")

$para1 = $d.Paragraphs.Item($para1Index)
$newPara2 = $para1.Range.InsertParagraphAfter()
$para2Index = $para1.Index + 1
$d.Paragraphs.Item($para2Index).Range.InsertBefore("As a user I want to be able to input a custom spell ")

# Place the _GoBack bookmark, collapsed, right at the end of the new text
# (immediately before the paragraph mark). Adding a bookmark at a collapsed
# range exactly one position before a paragraph's end is unreliable in this
# runtime, so we temporarily append a sentinel character, bookmark just
# before it, then remove the sentinel again.
$para2 = $d.Paragraphs.Item($para2Index)
$para2Range = $para2.Range
$sentinelPos = $para2Range.End - 1
$d.Range($sentinelPos, $sentinelPos).InsertAfter("~")

$para2Range2 = $d.Paragraphs.Item($para2Index).Range
$bmPos = $para2Range2.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$d.Range($bmPos, $bmPos + 1).Delete()
